$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = 112094770
$ws.Range("B2").Value = 94034
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 2869
$ws.Range("F2").Value = "Bollvitmossa"
$ws.Range("G2").Value = "Sphagnum wulfianum"
$ws.Range("H2").Value = "Girg."
$ws.Range("S2").Value = 50

# Row 3 updates
$ws.Range("A3").Value = 112094771
$ws.Range("B3").Value = 77636
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("S3").Value = 10

# Row 4 update
$ws.Range("B4").Value = 78725
